# VerifyReq_CopyRFQCancelledItems.xlsx edit
# Adds unitPrice / Quantity / UOMValue columns with sample data to the
# first worksheet (CopyRFQCancelledItems), mirroring the new RFQ cancelled
# item row copied into the cart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells
$ws.Range("C1").Value = "unitPrice"
$ws.Range("D1").Value = "Quantity"
$ws.Range("E1").Value = "UOMValue"

# New data row
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "CU-CUBIC"

# Size the new text columns to fit their contents
$ws.Range("C1:C2").EntireColumn.AutoFit()
$ws.Range("E1:E2").EntireColumn.AutoFit()

# Leave the selection on the cell just past the new data, as in the source
$ws.Range("F2").Select() | Out-Null
